# Applies the "added outline levels to styles" commit:
#   1. Re-points the lone body paragraph's style from ResumeContactLine to
#      ResumeWordJumble.
#   2. Adds an explicit OutlineLevel to a handful of the template's custom
#      paragraph styles (w:outlineLvl in the saved XML).
#   3. Nudges ResumePositionDescription's left indent from 630 -> 634 twips.

$d = $word.ActiveDocument

# -- 1. paragraph style swap -------------------------------------------------
$d.Paragraphs(1).Style = "ResumeWordJumble"

# -- 2. outline levels --------------------------------------------------------
# WdOutlineLevel is 1-based (wdOutlineLevel1 = 1 ... wdOutlineLevelBodyText = 10)
# and serialises to OOXML's zero-based <w:outlineLvl w:val="n"/>, so the COM
# value is always (target XML level + 1).
$outlineLevels = @{
    "ResumeSectionHeader"        = 1   # -> w:outlineLvl w:val="0"
    "ResumeWordJumble"           = 3   # -> w:outlineLvl w:val="2"
    "ResumeCompanyDescription"   = 4   # -> w:outlineLvl w:val="3"
    "ResumePositionDescription"  = 3   # -> w:outlineLvl w:val="2"
    "ResumeCompanyHeader"        = 2   # -> w:outlineLvl w:val="1"
    "ResumeContactLine"          = 2   # -> w:outlineLvl w:val="1"
    "ResumeJobAccomplishment"    = 3   # -> w:outlineLvl w:val="2"
    "ResumeName"                 = 1   # -> w:outlineLvl w:val="0"
}

foreach ($styleName in $outlineLevels.Keys) {
    $style = $d.Styles($styleName)
    $style.ParagraphFormat.OutlineLevel = $outlineLevels[$styleName]
}

# -- 3. ResumePositionDescription indent tweak (630 -> 634 twips) -----------
$posDescStyle = $d.Styles("ResumePositionDescription")
$posDescStyle.ParagraphFormat.LeftIndent = 634 / 20.0
